# Update the cryptos list with the latest scraped values (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (D value, E value). $null means "leave unchanged".
$rowUpdates = @{
    2  = @("30.491.71", "  +0.40%  ")
    3  = @("2.107.71",  "  +4.65%  ")
    4  = @($null,       "  +0.00%  ")
    5  = @("330.30",    $null)
    6  = @($null,       "  +0.03%  ")
    7  = @("0.5274",    "  +2.73%  ")
    8  = @("0.4401",    "  +3.29%  ")
    9  = @("0.08901",   "  +1.61%  ")
    10 = @("47.88",     "  +10.10%  ")
    11 = @("1.167",     "  +2.83%  ")
    12 = @("24.73",     "  +0.63%  ")
    13 = @("2.106.46",  "  +4.53%  ")
    14 = @("6.760",     "  +2.31%  ")
    15 = @("7.777",     "  +4.22%  ")
    16 = @("96.62",     "  +2.51%  ")
    17 = @($null,       "  +0.16%  ")
    18 = @("0.00001133","  +1.69%  ")
    19 = @("0.06639",   "  +1.72%  ")
    20 = @("19.08",     "  +0.95%  ")
    21 = @($null,       "  +0.10%  ")
    22 = @("6.321",     "  +1.74%  ")
    23 = @("30.550.36", "  +0.40%  ")
    24 = @("12.30",     "  +3.75%  ")
    25 = @("2.349",     "  +3.43%  ")
    26 = @("2.355.59",  "  +4.60%  ")
    27 = @("22.51",     "  +0.23%  ")
    28 = @("2.647",     "  +8.76%  ")
    29 = @("161.86",    "  -0.46%  ")
    30 = @("133.11",    "  +1.51%  ")
    31 = @("1.218",     "  +5.84%  ")
    32 = @($null,       "  +1.80%  ")
    33 = @("1.688",     "  +23.62%  ")
    34 = @("6.246",     "  +2.26%  ")
    35 = @("3.924",     "  +2.44%  ")
    36 = @("10.26",     "  +11.90%  ")
    37 = @("0.02587",   "  +2.14%  ")
    40 = @($null,       "  +2.83%  ")
    41 = @($null,       "  +2.98%  ")
    42 = @("0.6884",    "  +3.18%  ")
    43 = @("1.272",     "  +2.95%  ")
    44 = @($null,       "  +0.09%  ")
    45 = @("0.6421",    "  +3.84%  ")
    46 = @("14.05",     "  +2.79%  ")
    47 = @("2.220",     "  +1.01%  ")
    48 = @("3.631",     "  +0.02%  ")
    49 = @("1.256",     "  -0.22%  ")
    50 = @("1.217",     "  +10.08%  ")
    51 = @("82.53",     "  +1.59%  ")
}

foreach ($row in $rowUpdates.Keys) {
    $vals = $rowUpdates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        # Force text so strings like "330.30" or "47.88" aren't coerced into
        # numbers (which would drop the trailing zero / change formatting).
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
    }
    if ($null -ne $eVal) {
        # The padded "  +x.xx%  " strings never parse as numbers, so no
        # NumberFormat juggling is needed here.
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# Rows 38 and 39 swap places entirely (Hedera <-> InternetComputer(DFINITY))
# plus their Volume(1h) values change as well.
$ws.Cells.Item(38, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.513"
$ws.Cells.Item(38, 5).Value = "  +0.88%  "

$ws.Cells.Item(39, 2).Value = "Hedera"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06732"
$ws.Cells.Item(39, 5).Value = "  +1.03%  "
